$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6459696666666667
$ws.Range("H2").Value = 1.937909
$ws.Range("I2").Value = 0.2417008406295992
$ws.Range("J2").Value = 0.2417008406295992
$ws.Range("M2").Value = 0.232947
$ws.Range("N2").Value = 0.698841
$ws.Range("O2").Value = 0.2572219815457369
$ws.Range("P2").Value = 0.2572219815457369
$ws.Range("Q2").Value = 0.150476695941
$ws.Range("R2").Value = 1.354290263469
$ws.Range("S2").Value = 0.06217076916801585
$ws.Range("T2").Value = 0.06217076916801585
$ws.Range("G3").Value = 0.6459696666666667
$ws.Range("H3").Value = 1.937909
$ws.Range("I3").Value = 0.2417008406295992
$ws.Range("J3").Value = 0.2417008406295992
$ws.Range("M3").Value = 0.6726793333333333
$ws.Range("N3").Value = 2.018038
$ws.Range("O3").Value = 0.7427780184542632
$ws.Range("P3").Value = 0.7427780184542632
$ws.Range("Q3").Value = 0.4345304447268888
$ws.Range("R3").Value = 3.910774002541999
$ws.Range("S3").Value = 0.1795300714615833
$ws.Range("T3").Value = 0.1795300714615833
$ws.Range("G4").Value = 1.310047666666666
$ws.Range("I4").Value = 0.4901772306617775
$ws.Range("J4").Value = 0.4901772306617775
$ws.Range("M4").Value = 0.232947
$ws.Range("N4").Value = 0.698841
$ws.Range("O4").Value = 0.2572219815457369
$ws.Range("P4").Value = 0.2572219815457369
$ws.Range("Q4").Value = 0.305171673807
$ws.Range("R4").Value = 2.746545064263
$ws.Range("S4").Value = 0.1260843585794242
$ws.Range("T4").Value = 0.1260843585794242
$ws.Range("G5").Value = 1.310047666666666
$ws.Range("I5").Value = 0.4901772306617775
$ws.Range("J5").Value = 0.4901772306617775
$ws.Range("M5").Value = 0.6726793333333333
$ws.Range("N5").Value = 2.018038
$ws.Range("O5").Value = 0.7427780184542632
$ws.Range("P5").Value = 0.7427780184542632
$ws.Range("Q5").Value = 0.8812419910482221
$ws.Range("R5").Value = 7.931177919433998
$ws.Range("S5").Value = 0.3640928720823534
$ws.Range("T5").Value = 0.3640928720823534
$ws.Range("G6").Value = 0.4978893333333334
$ws.Range("H6").Value = 1.493668
$ws.Range("I6").Value = 0.1862939958592133
$ws.Range("J6").Value = 0.1862939958592133
$ws.Range("M6").Value = 0.232947
$ws.Range("N6").Value = 0.698841
$ws.Range("O6").Value = 0.2572219815457369
$ws.Range("P6").Value = 0.2572219815457369
$ws.Range("Q6").Value = 0.115981826532
$ws.Range("R6").Value = 1.043836438788
$ws.Range("S6").Value = 0.04791891076498014
$ws.Range("T6").Value = 0.04791891076498014
$ws.Range("G7").Value = 0.4978893333333334
$ws.Range("H7").Value = 1.493668
$ws.Range("I7").Value = 0.1862939958592133
$ws.Range("J7").Value = 0.1862939958592133
$ws.Range("M7").Value = 0.6726793333333333
$ws.Range("N7").Value = 2.018038
$ws.Range("O7").Value = 0.7427780184542632
$ws.Range("P7").Value = 0.7427780184542632
$ws.Range("Q7").Value = 0.3349198648204444
$ws.Range("R7").Value = 3.014278783383999
$ws.Range("S7").Value = 0.1383750850942331
$ws.Range("T7").Value = 0.1383750850942331
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2186933333333333
$ws.Range("H8").Value = 0.65608
$ws.Range("I8").Value = 0.08182793284941008
$ws.Range("J8").Value = 0.08182793284941006
$ws.Range("M8").Value = 0.232947
$ws.Range("N8").Value = 0.698841
$ws.Range("O8").Value = 0.2572219815457369
$ws.Range("P8").Value = 0.2572219815457369
$ws.Range("Q8").Value = 0.05094395592
$ws.Range("R8").Value = 0.4584956032800001
$ws.Range("S8").Value = 0.02104794303331676
$ws.Range("T8").Value = 0.02104794303331675
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2186933333333333
$ws.Range("H9").Value = 0.65608
$ws.Range("I9").Value = 0.08182793284941008
$ws.Range("J9").Value = 0.08182793284941006
$ws.Range("M9").Value = 0.6726793333333333
$ws.Range("N9").Value = 2.018038
$ws.Range("O9").Value = 0.7427780184542632
$ws.Range("P9").Value = 0.7427780184542632
$ws.Range("Q9").Value = 0.1471104856711111
$ws.Range("R9").Value = 1.32399437104
$ws.Range("S9").Value = 0.06077998981609332
$ws.Range("T9").Value = 0.06077998981609332
